# "did correction to user story data visualisation"
#
# 1. Age Groups sheet: add a new "0-25" age bracket row (1 guest) and keep
#    the bar chart's source range in sync with the extra row.
# 2. Nationalities sheet: correct the Switzerland guest count (6 -> 7) and
#    keep the pie chart's series formula pointing at the (unchanged) range.
# 3. Recurring sheet: replace the old "Recurring Guests" single-fact blurb
#    with a real recurring-guest visit table (name / visits / total nights
#    / last stay), including a YYYY-MM-DD date format on the last column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Age Groups
# ---------------------------------------------------------------------
$wsAge = $wb.Worksheets.Item("Age Groups")

$wsAge.Range("A4").Value = "0-25"
$wsAge.Range("B4").Value = 1

$ageChart = $wsAge.ChartObjects(1).Chart
$ageSeries = $ageChart.SeriesCollection(1)
$ageSeries.XValues = "='Age Groups'!`$A`$2:`$A`$4"
$ageSeries.Values = "='Age Groups'!`$B`$2:`$B`$4"

# ---------------------------------------------------------------------
# 2) Nationalities
# ---------------------------------------------------------------------
$wsNat = $wb.Worksheets.Item("Nationalities")

$wsNat.Range("B2").Value = 7

$natChart = $wsNat.ChartObjects(1).Chart
$natSeries = $natChart.SeriesCollection(1)
$natSeries.XValues = "=Nationalities!`$A`$2:`$A`$15"
$natSeries.Values = "=Nationalities!`$B`$2:`$B`$15"

# ---------------------------------------------------------------------
# 3) Recurring -> Guest visit table
# ---------------------------------------------------------------------
$wsRec = $wb.Worksheets.Item("Recurring")
$wsRec.Cells.Clear()

# Row 1: accessible table caption in A1 (unstyled) + real column headers.
$wsRec.Range("A1").Value = "Gastname, Anzahl Besuche, Gesamtanzahl Nächte, Letzter Aufenthalt"
$wsRec.Range("B1").Value = "Visits"
$wsRec.Range("C1").Value = "Total Nights"
$wsRec.Range("D1").Value = "Last Stay"

$recHeader = $wsRec.Range("B1:D1")
$recHeader.Font.Bold = $true
$recHeader.HorizontalAlignment = -4108
$recHeader.VerticalAlignment = -4160
$recHeader.Borders.LineStyle = 1

# Guest rows.
$wsRec.Range("A2").Value = "Muster Max"
$wsRec.Range("B2").Value = 2
$wsRec.Range("C2").Value = 10
$wsRec.Range("D2").Value = 45963

$wsRec.Range("A3").Value = "Schmidt Anna"
$wsRec.Range("B3").Value = 2
$wsRec.Range("C3").Value = 10
$wsRec.Range("D3").Value = 45964

$wsRec.Range("A4").Value = "Müller Thomas"
$wsRec.Range("B4").Value = 2
$wsRec.Range("C4").Value = 10
$wsRec.Range("D4").Value = 45965

$wsRec.Range("A5").Value = "Weber Laura"
$wsRec.Range("B5").Value = 2
$wsRec.Range("C5").Value = 10
$wsRec.Range("D5").Value = 45966

$wsRec.Range("A6").Value = "Fischer Michael"
$wsRec.Range("B6").Value = 2
$wsRec.Range("C6").Value = 11
$wsRec.Range("D6").Value = 45967

$wsRec.Range("A7").Value = "Nathan Jeremy"
$wsRec.Range("B7").Value = 2
$wsRec.Range("C7").Value = 11
$wsRec.Range("D7").Value = 45995

$wsRec.Range("D2:D7").NumberFormat = "YYYY-MM-DD"
